# Auto-generated edit script applying numeric market-price updates
# to the Halicarnassus_Profits leve-profit tables across all job sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 179.25
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 179.25
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 179.25
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -405.25
$ws.Range("H9").Value = 164.66667
$ws.Range("I9").Value = 79.333336
$ws.Range("J9").Value = 250
$ws.Range("K9").Value = 79.333336
$ws.Range("L9").Value = 250
$ws.Range("M9").Value = 89.666664
$ws.Range("N9").Value = -588
$ws.Range("H17").Value = 1015
$ws.Range("J17").Value = 1015
$ws.Range("L17").Value = 3045
$ws.Range("N17").Value = -3381
$ws.Range("H113").Value = 5512.5
$ws.Range("I113").Value = 5766.6665
$ws.Range("K113").Value = 5766.6665
$ws.Range("M113").Value = -2512.6665
$ws.Range("H125").Value = 932
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1176.9
$ws.Range("I2").Value = 1085.4445
$ws.Range("K2").Value = 1085.4445
$ws.Range("M2").Value = -972.4445000000001
$ws.Range("H61").Value = 1710.4445
$ws.Range("I61").Value = 1318.9333
$ws.Range("K61").Value = 1318.9333
$ws.Range("M61").Value = -1106.9333
$ws.Range("H74").Value = 2302.85
$ws.Range("I74").Value = 1566.5
$ws.Range("K74").Value = 1566.5
$ws.Range("M74").Value = -692.5
$ws.Range("H77").Value = 2302.85
$ws.Range("I77").Value = 1566.5
$ws.Range("K77").Value = 7832.5
$ws.Range("M77").Value = -3464.5
$ws.Range("H116").Value = 1176.9
$ws.Range("I116").Value = 1085.4445
$ws.Range("K116").Value = 1085.4445
$ws.Range("M116").Value = 1208.5555
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H123").Value = 1979899
$ws.Range("J123").Value = 1979899
$ws.Range("L123").Value = 1979899
$ws.Range("N123").Value = -1989699
$ws.Range("H132").Value = 1865.8334
$ws.Range("I132").Value = 1857.9412
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5573.8236
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3043.8236
$ws.Range("N132").Value = -11060
$ws.Range("H136").Value = 1710.4445
$ws.Range("I136").Value = 1318.9333
$ws.Range("K136").Value = 3956.7999
$ws.Range("M136").Value = -1406.7999
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1176.9
$ws.Range("I3").Value = 1085.4445
$ws.Range("K3").Value = 1085.4445
$ws.Range("M3").Value = -971.4445000000001
$ws.Range("H22").Value = 423.25
$ws.Range("I22").Value = 269.42856
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 269.42856
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -96.42856
$ws.Range("N22").Value = -1846
$ws.Range("H99").Value = 3286
$ws.Range("I99").Value = 3286
$ws.Range("K99").Value = 3286
$ws.Range("M99").Value = -1788
$ws.Range("H106").Value = 47037.6
$ws.Range("J106").Value = 47037.6
$ws.Range("L106").Value = 47037.6
$ws.Range("N106").Value = -49561.6
$ws.Range("H107").Value = 4998.9
$ws.Range("I107").Value = 1664.8334
$ws.Range("K107").Value = 1664.8334
$ws.Range("M107").Value = 255.1666
$ws.Range("H111").Value = 55000
$ws.Range("J111").Value = 55000
$ws.Range("L111").Value = 55000
$ws.Range("N111").Value = -63180
$ws.Range("H134").Value = 2195.75
$ws.Range("I134").Value = 942.6667
$ws.Range("K134").Value = 2828.0001
$ws.Range("M134").Value = -293.0001000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1573.1666
$ws.Range("I16").Value = 997.25
$ws.Range("K16").Value = 997.25
$ws.Range("M16").Value = -710.25
$ws.Range("H31").Value = 7234.875
$ws.Range("I31").Value = 2268.8
$ws.Range("K31").Value = 2268.8
$ws.Range("M31").Value = -1973.8
$ws.Range("H34").Value = 7234.875
$ws.Range("I34").Value = 2268.8
$ws.Range("K34").Value = 2268.8
$ws.Range("M34").Value = -2066.8
$ws.Range("H58").Value = 2477.0557
$ws.Range("I58").Value = 1246.8462
$ws.Range("K58").Value = 1246.8462
$ws.Range("M58").Value = -1043.8462
$ws.Range("H62").Value = 2900
$ws.Range("I62").Value = 2800
$ws.Range("K62").Value = 2800
$ws.Range("M62").Value = -2176
$ws.Range("H65").Value = 2900
$ws.Range("I65").Value = 2800
$ws.Range("K65").Value = 14000
$ws.Range("M65").Value = -10880
$ws.Range("H99").Value = 2388.25
$ws.Range("I99").Value = 2434.4546
$ws.Range("K99").Value = 2434.4546
$ws.Range("M99").Value = -936.4546
$ws.Range("H105").Value = 3004.4443
$ws.Range("I105").Value = 3004.4443
$ws.Range("K105").Value = 3004.4443
$ws.Range("M105").Value = -1257.4443
$ws.Range("H113").Value = 1573.1666
$ws.Range("I113").Value = 997.25
$ws.Range("K113").Value = 997.25
$ws.Range("M113").Value = 1172.75
$ws.Range("H122").Value = 931.2222
$ws.Range("I122").Value = 922.625
$ws.Range("K122").Value = 2767.875
$ws.Range("M122").Value = -317.875
$ws.Range("H126").Value = 2388.25
$ws.Range("I126").Value = 2434.4546
$ws.Range("K126").Value = 7303.3638
$ws.Range("M126").Value = -4833.3638
$ws.Range("H132").Value = 2000.75
$ws.Range("I132").Value = 2122.077
$ws.Range("K132").Value = 6366.231000000001
$ws.Range("M132").Value = -3836.231000000001
$ws.Range("H136").Value = 2477.0557
$ws.Range("I136").Value = 1246.8462
$ws.Range("K136").Value = 3740.5386
$ws.Range("M136").Value = -1190.5386

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2884.5715
$ws.Range("I113").Value = 1644.6666
$ws.Range("K113").Value = 1644.6666
$ws.Range("M113").Value = 525.3334
$ws.Range("H132").Value = 2168.2
$ws.Range("I132").Value = 1909.2222
$ws.Range("K132").Value = 5727.6666
$ws.Range("M132").Value = -3197.6666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8237.6
$ws.Range("I7").Value = 7922.25
$ws.Range("J7").Value = 9499
$ws.Range("K7").Value = 7922.25
$ws.Range("L7").Value = 9499
$ws.Range("M7").Value = -7810.25
$ws.Range("N7").Value = -9723
$ws.Range("H22").Value = 2800
$ws.Range("J22").Value = 3500
$ws.Range("L22").Value = 3500
$ws.Range("N22").Value = -4090
$ws.Range("H27").Value = 2800
$ws.Range("J27").Value = 3500
$ws.Range("L27").Value = 3500
$ws.Range("N27").Value = -3714
$ws.Range("H40").Value = 6665.1665
$ws.Range("I40").Value = 4996
$ws.Range("K40").Value = 4996
$ws.Range("M40").Value = -4860
$ws.Range("H61").Value = 3969.0908
$ws.Range("I61").Value = 2457.5
$ws.Range("K61").Value = 2457.5
$ws.Range("M61").Value = -2255.5
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H100").Value = 6599.8887
$ws.Range("I100").Value = 2649.6667
$ws.Range("K100").Value = 2649.6667
$ws.Range("M100").Value = -2108.6667
$ws.Range("H113").Value = 3969.0908
$ws.Range("I113").Value = 2457.5
$ws.Range("K113").Value = 2457.5
$ws.Range("M113").Value = -287.5
$ws.Range("H122").Value = 2746
$ws.Range("I122").Value = 2688
$ws.Range("K122").Value = 8064
$ws.Range("M122").Value = -5614
$ws.Range("H126").Value = 8237.6
$ws.Range("I126").Value = 7922.25
$ws.Range("J126").Value = 9499
$ws.Range("K126").Value = 23766.75
$ws.Range("L126").Value = 28497
$ws.Range("M126").Value = -21296.75
$ws.Range("N126").Value = -33437
$ws.Range("H132").Value = 4101
$ws.Range("I132").Value = 4015
$ws.Range("K132").Value = 12045
$ws.Range("M132").Value = -9515

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 7500
$ws.Range("J49").Value = 7500
$ws.Range("L49").Value = 7500
$ws.Range("N49").Value = -7960
$ws.Range("H126").Value = 6912.0386
$ws.Range("I126").Value = 6086
$ws.Range("J126").Value = 7738.077
$ws.Range("K126").Value = 18258
$ws.Range("L126").Value = 23214.231
$ws.Range("M126").Value = -15788
$ws.Range("N126").Value = -28154.231
$ws.Range("H132").Value = 1789.8846
$ws.Range("I132").Value = 1701.48
$ws.Range("K132").Value = 5104.440000000001
$ws.Range("M132").Value = -2574.440000000001

